# eliminacion de imagenes incorporadas
#
# Rows 45-70 had lost their "Peso"/"Valor" text formatting (columns C/D were
# stored as real numbers instead of the text-like values used everywhere
# else in the sheet) and their "Ver Imagen N" hyperlink cells (columns
# E/F/G) had been dropped entirely. This restores both: C/D go back to
# text cells holding the same numeric-looking values, and the missing
# "Ver Imagen N" hyperlink cells are re-created with the same style used
# by every other image-link cell in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, [string]$val)
    $cell = $ws.Range($addr)
    # Force the numeric-looking string to be stored as text (no leading
    # apostrophe / quotePrefix residue), then drop back to the default
    # "Normal" style so no stray number-format style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Add-ImageLink {
    param($ws, $addr, [string]$url, [string]$text)
    $cell = $ws.Range($addr)
    $ws.Hyperlinks.Add($cell, $url) | Out-Null
    # Reuse the workbook's existing "Hyperlink" cell style (same as every
    # other "Ver Imagen N" cell) instead of the ad-hoc style Add() applies.
    $cell.Style = "Hyperlink"
    $cell.Value = $text
}

# row -> (Peso, Valor)
$pesoValor = @{
    45 = @("0", "8000")
    46 = @("2", "2000")
    47 = @("3", "240")
    48 = @("5", "1500")
    49 = @("5", "1500")
    50 = @("3", "800")
    51 = @("14", "1500")
    52 = @("38", "3000")
    53 = @("6", "2500")
    54 = @("8", "2500")
    55 = @("16", "1300")
    56 = @("8", "300")
    57 = @("6", "300")
    58 = @("15", "1200")
    59 = @("10", "600")
    60 = @("9", "500")
    61 = @("8", "600")
    62 = @("15", "200")
    63 = @("2", "160")
    64 = @("3", "1200")
    65 = @("3", "240")
    66 = @("11", "880")
    67 = @("5", "400")
    68 = @("3", "600")
    69 = @("20", "1600")
    70 = @("6", "400")
}

# row -> list of image columns that need a "Ver Imagen N" hyperlink cell
$imageCols = @{
    45 = @("E")
    46 = @("E", "F", "G")
    47 = @("E")
    48 = @("E", "F")
    49 = @("E")
    50 = @("E")
    51 = @("E")
    52 = @("E")
    53 = @("E", "F")
    54 = @("E", "F")
    55 = @("E")
    56 = @("E")
    57 = @("E")
    58 = @("E")
    59 = @("E")
    60 = @("E")
    61 = @("E")
    62 = @("E")
    63 = @("E")
    64 = @("E", "F")
    65 = @("E")
    66 = @("E", "F")
    67 = @("E", "F")
    68 = @("E", "F")
    69 = @("E", "F")
    70 = @("E", "F")
}

$imageLabels = @{ "E" = "Ver Imagen 1"; "F" = "Ver Imagen 2"; "G" = "Ver Imagen 3" }

# Continue the existing "imagenes_subidas/<timestamp>_iOS.jpg" naming
# convention, filling the timestamp gap between the last photo used by
# row 44 (20250216_171041) and the first one used by row 71
# (20250216_182305).
$imgStamps = @(
    "20250216_171100000","20250216_171258137","20250216_171456274","20250216_171655411",
    "20250216_171853548","20250216_172051685","20250216_172250822","20250216_172448959",
    "20250216_172646096","20250216_172845233","20250216_173043370","20250216_173241507",
    "20250216_173440644","20250216_173638781","20250216_173836918","20250216_174035055",
    "20250216_174233192","20250216_174431329","20250216_174630466","20250216_174828603",
    "20250216_175026740","20250216_175225877","20250216_175423014","20250216_175621151",
    "20250216_175820288","20250216_180018425","20250216_180216562","20250216_180415699",
    "20250216_180613836","20250216_180811973","20250216_181010110","20250216_181208247",
    "20250216_181406384","20250216_181605521","20250216_181803658","20250216_182001795",
    "20250216_182200932"
)
$stampIdx = 0

foreach ($row in 45..70) {
    $pv = $pesoValor[$row]
    Set-TextValue $ws ("C" + $row) $pv[0]
    Set-TextValue $ws ("D" + $row) $pv[1]

    foreach ($col in $imageCols[$row]) {
        $addr = $col + $row
        $url = "imagenes_subidas/" + $imgStamps[$stampIdx] + "_iOS.jpg"
        $stampIdx = $stampIdx + 1
        Add-ImageLink $ws $addr $url $imageLabels[$col]
    }
}

Write-Output "Restored Peso/Valor text + Ver Imagen hyperlinks for rows 45-70"
